$wb = $excel.ActiveWorkbook

# Add the new "Description key" sheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Description key"

# Header row (bold, mirrors the style used on the other key sheets)
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Description"
$ws.Range("A1:C1").Font.Bold = $true

# Data rows
$ws.Range("A2").Value = "Northbound"
$ws.Range("C3").Value = "Traffic toward Lagos Island"
$ws.Range("A3").Value = "Southbound"
$ws.Range("C2").Value = "Traffic to Berger, Ikeja (Mainland)"
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2

$ws.Range("C3").Select()

# Make the newly added sheet the active/selected tab
$ws.Activate()
